$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The paragraph that reads:
#   "The are most suitable would be between localities ..."
# becomes:
#   "The most suitable location would be between localities ..."
# and the hidden "_GoBack" bookmark (which used to sit, empty, in the
# following Heading1 paragraph) moves to sit right after the newly
# inserted word " location" in this paragraph - i.e. it now marks the
# position of the user's last edit.
# ---------------------------------------------------------------------------

$oldOpening = "The are most suitable would be between localities"
$newOpening = "The most suitable location would be between localities"

$find = $d.Content
$found = $find.Find.Execute($oldOpening, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newOpening, 2)

if ($found) {
    # Re-locate the (now updated) sentence so we know exactly where it
    # starts, then work out character offsets relative to that start.
    $sentence = $d.Content
    $sentence.Find.Execute($newOpening, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
    $start = $sentence.Start

    $part1 = "The most suitable"          # run 1 (unchanged wording)
    $part2 = " location"                  # run 2 (newly typed word)
    $part3 = " would be "                 # run 3 (rest of the original run)

    $afterPart1 = $start + $part1.Length
    $afterPart2 = $afterPart1 + $part2.Length
    $afterPart3 = $afterPart2 + $part3.Length

    # Move the hidden "_GoBack" bookmark to sit between " location" and
    # " would be " - this is also what forces the XML run to split at
    # that point, matching real Word's behaviour when a bookmark sits in
    # the middle of typed text.
    $goBackRange = $d.Range($afterPart2, $afterPart2)
    $d.Bookmarks.Add("_GoBack", $goBackRange)

    # Force the remaining two run boundaries (between "The most suitable"
    # and " location", and between " would be " and "between localities")
    # using a throw-away bookmark that is immediately removed again - it
    # leaves the run split behind without leaving any bookmark markup.
    $split1 = $d.Range($afterPart1, $afterPart1)
    $d.Bookmarks.Add("ZZZtempSplit1", $split1)
    $d.Bookmarks("ZZZtempSplit1").Delete()

    $split2 = $d.Range($afterPart3, $afterPart3)
    $d.Bookmarks.Add("ZZZtempSplit2", $split2)
    $d.Bookmarks("ZZZtempSplit2").Delete()
}
